# Update cryptocurrency price (column D) and Volume(1h) (column E) values
# for the "cryptos" symbol-list refresh, as produced by the GitHub Actions
# scheduled job on Fri Jan 27 17:25:26 UTC 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = "307.13";      E = "0.91%" }
    3  = @{ D = "36.26";       E = "1.18%" }
    4  = @{ D = "5.060";       E = "-0.16%" }
    5  = @{ D = "0.08079";     E = "0.23%" }
    6  = @{ D = "2.015";       E = "4.22%" }
    7  = @{ D = "7.828";       E = "-0.11%" }
    8  = @{ D = "0.9264";      E = "-0.28%" }
    9  = @{ D = "0.1465";      E = "11.99%" }
    10 = @{ D = "0.1923";      E = "0.93%" }
    11 = @{ D = "0.09080";     E = "-1.49%" }
    12 = @{ D = "0.03438";     E = "-1.19%" }
    13 = @{ D = "0.09917";     E = "0.18%" }
    14 = @{ D = "0.001409";    E = "-0.59%" }
    15 = @{ D = "0.006202";    E = "-6.62%" }
    16 = @{ D = "3.840";       E = "6.34%" }
    17 = @{ D = "4.143";       E = "-0.29%" }
    18 = @{ E = "11.60%" }
    19 = @{ D = "0.3451";      E = "0.81%" }
    20 = @{ D = "0.1335";      E = "2.39%" }
    21 = @{ D = "4.793";       E = "-7.28%" }
    22 = @{ D = "0.2338";      E = "-7.68%" }
    23 = @{ D = "0.04359";     E = "-1.28%" }
    24 = @{ D = "0.001228";    E = "-0.52%" }
    25 = @{ D = "0.004302";    E = "-8.45%" }
    27 = @{ D = "0.0001299";   E = "-0.25%" }
    39 = @{ D = "0.02018";     E = "1.11%" }
    40 = @{ D = "0.05147";     E = "-0.74%" }
    41 = @{ D = "0.007492";    E = "-1.61%" }
    42 = @{ E = "-0.15%" }
    43 = @{ D = "0.1363";      E = "-0.05%" }
    44 = @{ D = "0.002148";    E = "2.13%" }
    45 = @{ D = "0.009930";    E = "-7.51%" }
    46 = @{ D = "0.00006272";  E = "-0.16%" }
    47 = @{ D = "0.00000000749"; E = "-0.19%" }
    49 = @{ D = "0.001249";    E = "-22.02%" }
    50 = @{ D = "0.00002098";  E = "-0.19%" }
    51 = @{ D = "0.0001998";   E = "-0.19%" }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $cell = $ws.Range("$col$row")
        # Force text storage so numeric-looking strings (prices, percents)
        # stay literal text instead of being reinterpreted as numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col]
    }
}
